$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab14")

# --- Update header labels: "most recent measure 2012-21" -> "...2013-22" ---
$ws.Range("C2").Value = "Cantril life ladder, most recent measure 2013-22"
$ws.Range("D2").Value = "Negative affect, most recent measure 2013-22"
$ws.Range("E2").Value = "Positive affect, most recent measure 2013-22"
$ws.Range("F2").Value = "Freedom to make life choices, most recent measure 2013-22"
$ws.Range("G2").Value = "Generosity, most recent measure 2013-22"
$ws.Range("H2").Value = "Social support, most recent measure 2013-22"

# --- Fix footnote text (regional economic communities) ---
$ws.Range("A103").Value = "Regional Economic Communities:CEN-SAD = `"Community of Sahel-Saharan States`";COMESA = `"Common Market for Eastern and Southern Africa`";EAC = `"East African Community`";ECCAS = `"Economic Community of Central African States`";ECOWAS = `"Economic Community of West African States`";IGAD = `"Intergovernmental Authority on Development`";SADC = `"Southern African Development Community`";UMA = `"Arab Maghreb Union`";PALOP = `"Pa&gt;ses Africanos de L&gt;ngua Oficial Portuguesa`";ASEAN = `"Association of Southeast Asian Nations`";MERCOSUR = `"Mercado Com&gt;n del Sur`".EU27 = `"European Union (27 members)`".OECD = `"Organisation for Economic Co-operation and Development`"."

# --- Updated statistical values (refreshed source data) ---
$ws.Range("C63").Value = 5.8902884258164301
$ws.Range("D63").Value = 0.29460977791053
$ws.Range("E63").Value = 0.72038237356915003
$ws.Range("F63").Value = 0.82672665553671998
$ws.Range("G63").Value = [double]"-5.6752783994999997E-3"
$ws.Range("H63").Value = 0.85126551735067002

$ws.Range("C64").Value = 5.7897235371849796
$ws.Range("D64").Value = 0.31858937713232999
$ws.Range("E64").Value = 0.79291529005224004
$ws.Range("F64").Value = 0.83787702159447996
$ws.Range("G64").Value = [double]"-5.2704396911599997E-2"
$ws.Range("H64").Value = 0.82073891975663005

$ws.Range("C66").Value = 5.4690353593518699
$ws.Range("D66").Value = 0.31112098031573998
$ws.Range("E66").Value = 0.70601781696276
$ws.Range("F66").Value = 0.78623798934669997
$ws.Range("G66").Value = [double]"-1.6922179279000001E-3"
$ws.Range("H66").Value = 0.80125126571623995

$ws.Range("C77").Value = 5.7319341659545904
$ws.Range("D77").Value = 0.35167565345764001
$ws.Range("E77").Value = 0.78529569506645003
$ws.Range("F77").Value = 0.81471749544144001
$ws.Range("G77").Value = [double]"-7.9204652272199999E-2"
$ws.Range("H77").Value = 0.84500606060028005

$ws.Range("C81").Value = 5.6538367271423402
$ws.Range("D81").Value = 0.28696681279689001
$ws.Range("E81").Value = 0.69262163713573999
$ws.Range("F81").Value = 0.80691246315836995
$ws.Range("G81").Value = [double]"6.7111144307999996E-3"
$ws.Range("H81").Value = 0.86383996903895999

$ws.Range("C89").Value = 5.5868839945111999
$ws.Range("D89").Value = 0.30767310454565
$ws.Range("E89").Value = 0.69209085667835002
$ws.Range("F89").Value = 0.80759759034429002
$ws.Range("G89").Value = [double]"-5.11234869116E-2"
$ws.Range("H89").Value = 0.84280766929898998

$ws.Range("D90").Value = 0.25534762272780998
$ws.Range("E90").Value = 0.74723589420319003
$ws.Range("F90").Value = 0.86335812644525001
$ws.Range("G90").Value = [double]"-2.9786770665899999E-2"
$ws.Range("H90").Value = 0.91762280464171997

$ws.Range("C94").Value = 5.6236491501331303
$ws.Range("D94").Value = 0.26179147192410002
$ws.Range("E94").Value = 0.74055125032152003
$ws.Range("F94").Value = 0.84844066202641
$ws.Range("G94").Value = [double]"4.7176104970279999E-2"
$ws.Range("H94").Value = 0.82260600477457002
